$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that disappear entirely ("empty" at row 10 and "empty; " at row 37),
# shifting the remaining rows up to close the gaps.
$ws.Rows.Item(37).Delete()
$ws.Rows.Item(10).Delete()

# Final data for A2:B36 (problem/DOI header in row 1 is unchanged)
$data = @(
    @("no abstract in article", 54),
    @("includes multiple languages", 19),
    @("article includes abstract", 18),
    @("front matter", 13),
    @("MathML", 8),
    @("bad url", 8),
    @("article version includes two languages", 7),
    @("abstract not included in PDF", 6),
    @("includes hyperlinks", 5),
    @("includes disclosures; includes figure", 5),
    @("includes disclosures", 5),
    @("includes HTML", 4),
    @("includes full citation format in the metadata", 4),
    @("article title used as abstract", 4),
    @("can't verify", 4),
    @("includes U+2009 Thin Space General Punc char", 3),
    @("first element is empty, but others have abstract", 3),
    @("U+2009 thin space from general punc; ", 3),
    @("includes references", 2),
    @("abstract is longer on the landing page", 2),
    @("includes funding acknowledgement", 2),
    @("hiragana, katakana, CJK Ideograph, latin", 1),
    @("section title used as abstract", 1),
    @("marked as retracted", 1),
    @("includes tex-math, not MathML", 1),
    @("includes tables; includes full citation format in the metadata", 1),
    @("includes table", 1),
    @("includes hyperlinks; includes funding acknowledgement; includes figures; ", 1),
    @("U+202F Zero break space from general punc", 1),
    @("MathML; includes copyright statement", 1),
    @("abstract is longer in the metadata", 1),
    @("includes figure, acknowledgements", 1),
    @("extra \n char render in REST API, but not downloaded data", 1),
    @("includes U+2771 Dingbat char", 1),
    @("test doi", 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
